# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the "Periodo Mora" / "Valor Mora" rows (16-22) ---
# The period column (E) is reversed (2106..2112 -> 2112..2106) and the
# matching "Valor Mora" amount follows the same row it was already on
# (i.e. the whole E/F pair per row is mirrored top<->bottom).
$firstRow = 16
$lastRow = 22

$periods = @()
$amounts = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += ,($ws.Range("E$r").Value())
    $amounts += ,($ws.Range("F$r").Value())
}

$count = $periods.Length
for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $mirrorIndex = $count - 1 - $i
    $ws.Range("E$r").Value = $periods[$mirrorIndex]
    $ws.Range("F$r").Value = $amounts[$mirrorIndex]
}

# --- Adjust column widths slightly (data/layout refresh) ---
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
